# Apply updated Kraken market-price snapshot values to the per-class Leve profit sheets.
# Generated from the authoritative cell-level diff (sheet, row, column -> new value).
$wb = $excel.ActiveWorkbook


# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Cells.Item(15, 8).Value = 342.3846   # H15
$ws.Cells.Item(15, 9).Value = 342.3846   # I15
$ws.Cells.Item(15, 11).Value = 1027.1538   # K15
$ws.Cells.Item(15, 13).Value = -858.1538   # M15
# Row 28
$ws.Cells.Item(28, 8).Value = 35.6   # H28
$ws.Cells.Item(28, 9).Value = 35.6   # I28
$ws.Cells.Item(28, 10).Value = 0   # J28
$ws.Cells.Item(28, 11).Value = 35.6   # K28
$ws.Cells.Item(28, 12).Value = 0   # L28
$ws.Cells.Item(28, 13).ClearContents()   # M28 (cell removed)
$ws.Cells.Item(28, 14).Value = 449.4   # N28
# Row 43
$ws.Cells.Item(43, 8).Value = 3000   # H43
$ws.Cells.Item(43, 10).Value = 3000   # J43
$ws.Cells.Item(43, 12).Value = 3000   # L43
$ws.Cells.Item(43, 14).Value = -3138   # N43
# Row 53
$ws.Cells.Item(53, 8).Value = 509.83334   # H53
$ws.Cells.Item(53, 9).Value = 509.83334   # I53
$ws.Cells.Item(53, 11).Value = 509.83334   # K53
$ws.Cells.Item(53, 13).Value = 127.16666   # M53
# Row 54
$ws.Cells.Item(54, 8).Value = 50000   # H54
$ws.Cells.Item(54, 10).Value = 50000   # J54
$ws.Cells.Item(54, 12).Value = 50000   # L54
$ws.Cells.Item(54, 14).Value = -50972   # N54
# Row 97
$ws.Cells.Item(97, 8).Value = 7253.1665   # H97
$ws.Cells.Item(97, 10).Value = 7253.1665   # J97
$ws.Cells.Item(97, 12).Value = 21759.4995   # L97
$ws.Cells.Item(97, 14).Value = -22751.4995   # N97
# Row 106
$ws.Cells.Item(106, 8).Value = 1990   # H106
$ws.Cells.Item(106, 9).Value = 1980   # I106
$ws.Cells.Item(106, 11).Value = 1980   # K106
$ws.Cells.Item(106, 13).Value = -1349   # M106
# Row 112
$ws.Cells.Item(112, 8).Value = 3858   # H112
$ws.Cells.Item(112, 10).Value = 3858   # J112
$ws.Cells.Item(112, 12).Value = 11574   # L112
$ws.Cells.Item(112, 14).Value = -13790   # N112
# Row 116
$ws.Cells.Item(116, 8).Value = 4966.6665   # H116
$ws.Cells.Item(116, 10).Value = 4900   # J116
$ws.Cells.Item(116, 12).Value = 4900   # L116
$ws.Cells.Item(116, 14).Value = -11784   # N116
# Row 132
$ws.Cells.Item(132, 8).Value = 8628.143   # H132
$ws.Cells.Item(132, 9).Value = 6224.25   # I132
$ws.Cells.Item(132, 11).Value = 18672.75   # K132
$ws.Cells.Item(132, 13).Value = -16142.75   # M132
# Row 135
$ws.Cells.Item(135, 8).Value = 1517.75   # H135
$ws.Cells.Item(135, 9).Value = 1214.375   # I135
$ws.Cells.Item(135, 11).Value = 10929.375   # K135
$ws.Cells.Item(135, 13).Value = -8394.375   # M135
# Row 137
$ws.Cells.Item(137, 8).Value = 3672.2068   # H137
$ws.Cells.Item(137, 10).Value = 3731.5264   # J137
$ws.Cells.Item(137, 12).Value = 11194.5792   # L137
$ws.Cells.Item(137, 14).Value = -16294.5792   # N137

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 47
$ws.Cells.Item(47, 8).Value = 0   # H47
$ws.Cells.Item(47, 10).Value = 0   # J47
$ws.Cells.Item(47, 12).ClearContents()   # L47 (cell removed)
$ws.Cells.Item(47, 14).Value = 0   # N47

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 19405.268   # H31
$ws.Cells.Item(31, 9).Value = 10957.083   # I31
$ws.Cells.Item(31, 11).Value = 10957.083   # K31
$ws.Cells.Item(31, 13).Value = -10662.083   # M31
# Row 34
$ws.Cells.Item(34, 8).Value = 19405.268   # H34
$ws.Cells.Item(34, 9).Value = 10957.083   # I34
$ws.Cells.Item(34, 11).Value = 10957.083   # K34
$ws.Cells.Item(34, 13).Value = -10755.083   # M34
# Row 58
$ws.Cells.Item(58, 8).Value = 5000   # H58
$ws.Cells.Item(58, 9).Value = 5000   # I58
$ws.Cells.Item(58, 11).Value = 5000   # K58
$ws.Cells.Item(58, 13).Value = -4797   # M58
# Row 94
$ws.Cells.Item(94, 8).Value = 1999.5   # H94
$ws.Cells.Item(94, 10).Value = 1999.5   # J94
$ws.Cells.Item(94, 12).Value = 1999.5   # L94
$ws.Cells.Item(94, 14).Value = -2901.5   # N94
# Row 122
$ws.Cells.Item(122, 8).Value = 2245.5   # H122
$ws.Cells.Item(122, 9).Value = 2245.5   # I122
$ws.Cells.Item(122, 11).Value = 6736.5   # K122
$ws.Cells.Item(122, 13).Value = -4286.5   # M122
# Row 132
$ws.Cells.Item(132, 8).Value = 2828.0908   # H132
$ws.Cells.Item(132, 9).Value = 2893.625   # I132
$ws.Cells.Item(132, 10).Value = 2653.3333   # J132
$ws.Cells.Item(132, 11).Value = 8680.875   # K132
$ws.Cells.Item(132, 12).Value = 7959.999899999999   # L132
$ws.Cells.Item(132, 13).Value = -6150.875   # M132
$ws.Cells.Item(132, 14).Value = -13019.9999   # N132
# Row 136
$ws.Cells.Item(136, 8).Value = 5000   # H136
$ws.Cells.Item(136, 9).Value = 5000   # I136
$ws.Cells.Item(136, 11).Value = 15000   # K136
$ws.Cells.Item(136, 13).Value = -12450   # M136

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Cells.Item(5, 8).Value = 5000   # H5
$ws.Cells.Item(5, 9).Value = 5000   # I5
$ws.Cells.Item(5, 10).Value = 0   # J5
$ws.Cells.Item(5, 11).Value = 15000   # K5
$ws.Cells.Item(5, 12).Value = 0   # L5
$ws.Cells.Item(5, 13).ClearContents()   # M5 (cell removed)
$ws.Cells.Item(5, 14).Value = -14888   # N5
# Row 12
$ws.Cells.Item(12, 8).Value = 49.5   # H12
$ws.Cells.Item(12, 10).Value = 49.5   # J12
$ws.Cells.Item(12, 12).Value = 148.5   # L12
$ws.Cells.Item(12, 14).Value = -494.5   # N12
# Row 80
$ws.Cells.Item(80, 8).Value = 5750.625   # H80
$ws.Cells.Item(80, 9).Value = 0   # I80
$ws.Cells.Item(80, 10).Value = 5750.625   # J80
$ws.Cells.Item(80, 11).Value = 0   # K80
$ws.Cells.Item(80, 12).ClearContents()   # L80 (cell removed)
$ws.Cells.Item(80, 13).Value = 17251.875   # M80
$ws.Cells.Item(80, 14).Value = -19123.875   # N80
# Row 83
$ws.Cells.Item(83, 8).Value = 5750.625   # H83
$ws.Cells.Item(83, 9).Value = 0   # I83
$ws.Cells.Item(83, 10).Value = 5750.625   # J83
$ws.Cells.Item(83, 11).Value = 0   # K83
$ws.Cells.Item(83, 12).ClearContents()   # L83 (cell removed)
$ws.Cells.Item(83, 13).Value = 51755.625   # M83
$ws.Cells.Item(83, 14).Value = -61115.625   # N83
# Row 113
$ws.Cells.Item(113, 8).Value = 1666.3334   # H113
$ws.Cells.Item(113, 9).Value = 1000   # I113
$ws.Cells.Item(113, 10).Value = 1999.5   # J113
$ws.Cells.Item(113, 11).Value = 3000   # K113
$ws.Cells.Item(113, 12).Value = 5998.5   # L113
$ws.Cells.Item(113, 13).Value = -830   # M113
$ws.Cells.Item(113, 14).Value = -10338.5   # N113
# Row 122
$ws.Cells.Item(122, 8).Value = 8224.75   # H122
$ws.Cells.Item(122, 9).Value = 8000   # I122
$ws.Cells.Item(122, 10).Value = 8449.5   # J122
$ws.Cells.Item(122, 11).Value = 72000   # K122
$ws.Cells.Item(122, 12).Value = 76045.5   # L122
$ws.Cells.Item(122, 13).Value = -69550   # M122
$ws.Cells.Item(122, 14).Value = -80945.5   # N122
# Row 135
$ws.Cells.Item(135, 8).Value = 5000   # H135
$ws.Cells.Item(135, 9).Value = 5000   # I135
$ws.Cells.Item(135, 10).Value = 0   # J135
$ws.Cells.Item(135, 11).Value = 45000   # K135
$ws.Cells.Item(135, 12).Value = 0   # L135
$ws.Cells.Item(135, 13).ClearContents()   # M135 (cell removed)
$ws.Cells.Item(135, 14).Value = -42465   # N135
# Row 139
$ws.Cells.Item(139, 8).Value = 3999   # H139
$ws.Cells.Item(139, 9).Value = 0   # I139
$ws.Cells.Item(139, 11).Value = 0   # K139
$ws.Cells.Item(139, 13).ClearContents()   # M139 (cell removed)

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Cells.Item(102, 8).Value = 4132.6665   # H102
$ws.Cells.Item(102, 9).Value = 4132.6665   # I102
$ws.Cells.Item(102, 11).Value = 4132.6665   # K102
$ws.Cells.Item(102, 13).Value = -2510.6665   # M102
# Row 122
$ws.Cells.Item(122, 8).Value = 8111.684   # H122
$ws.Cells.Item(122, 10).Value = 35283   # J122
$ws.Cells.Item(122, 12).Value = 105849   # L122
$ws.Cells.Item(122, 14).Value = -110749   # N122
# Row 132
$ws.Cells.Item(132, 8).Value = 0   # H132
$ws.Cells.Item(132, 9).Value = 0   # I132
$ws.Cells.Item(132, 11).Value = 0   # K132
$ws.Cells.Item(132, 13).ClearContents()   # M132 (cell removed)

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Cells.Item(7, 8).Value = 5086.8   # H7
$ws.Cells.Item(7, 9).Value = 4845   # I7
$ws.Cells.Item(7, 11).Value = 4845   # K7
$ws.Cells.Item(7, 13).Value = -4733   # M7
# Row 16
$ws.Cells.Item(16, 8).Value = 1185.6   # H16
$ws.Cells.Item(16, 9).Value = 1315.5   # I16
$ws.Cells.Item(16, 10).Value = 666   # J16
$ws.Cells.Item(16, 11).Value = 1315.5   # K16
$ws.Cells.Item(16, 12).Value = 666   # L16
$ws.Cells.Item(16, 13).Value = -1145.5   # M16
$ws.Cells.Item(16, 14).Value = -1006   # N16
# Row 40
$ws.Cells.Item(40, 8).Value = 6552   # H40
$ws.Cells.Item(40, 9).Value = 6138.2856   # I40
$ws.Cells.Item(40, 11).Value = 6138.2856   # K40
$ws.Cells.Item(40, 13).Value = -6002.2856   # M40
# Row 46
$ws.Cells.Item(46, 8).Value = 906.8333   # H46
$ws.Cells.Item(46, 9).Value = 814   # I46
$ws.Cells.Item(46, 10).Value = 999.6667   # J46
$ws.Cells.Item(46, 11).Value = 814   # K46
$ws.Cells.Item(46, 12).Value = 999.6667   # L46
$ws.Cells.Item(46, 13).Value = -626   # M46
$ws.Cells.Item(46, 14).Value = -1375.6667   # N46
# Row 122
$ws.Cells.Item(122, 8).Value = 6878.6   # H122
$ws.Cells.Item(122, 10).Value = 6833.3335   # J122
$ws.Cells.Item(122, 12).Value = 20500.0005   # L122
$ws.Cells.Item(122, 14).Value = -25400.0005   # N122
# Row 126
$ws.Cells.Item(126, 8).Value = 5086.8   # H126
$ws.Cells.Item(126, 9).Value = 4845   # I126
$ws.Cells.Item(126, 11).Value = 14535   # K126
$ws.Cells.Item(126, 13).Value = -12065   # M126
# Row 132
$ws.Cells.Item(132, 8).Value = 37833.168   # H132
$ws.Cells.Item(132, 9).Value = 35399.8   # I132
$ws.Cells.Item(132, 11).Value = 106199.4   # K132
$ws.Cells.Item(132, 13).Value = -103669.4   # M132

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 2
$ws.Cells.Item(2, 8).Value = 55000   # H2
$ws.Cells.Item(2, 9).Value = 55000   # I2
$ws.Cells.Item(2, 11).Value = 55000   # K2
$ws.Cells.Item(2, 13).Value = -54888   # M2
# Row 4
$ws.Cells.Item(4, 8).Value = 553.75   # H4
$ws.Cells.Item(4, 9).Value = 61.42857   # I4
$ws.Cells.Item(4, 11).Value = 61.42857   # K4
$ws.Cells.Item(4, 13).Value = 51.57143   # M4
# Row 54
$ws.Cells.Item(54, 8).Value = 21330   # H54
$ws.Cells.Item(54, 9).Value = 825   # I54
$ws.Cells.Item(54, 11).Value = 825   # K54
$ws.Cells.Item(54, 13).Value = -305   # M54
# Row 122
$ws.Cells.Item(122, 8).Value = 501371.5   # H122
$ws.Cells.Item(122, 9).Value = 667162.3   # I122
$ws.Cells.Item(122, 11).Value = 2001486.9   # K122
$ws.Cells.Item(122, 13).Value = -1999036.9   # M122
# Row 136
$ws.Cells.Item(136, 8).Value = 3060.9167   # H136
$ws.Cells.Item(136, 9).Value = 2673.2   # I136
$ws.Cells.Item(136, 11).Value = 8019.599999999999   # K136
$ws.Cells.Item(136, 13).Value = -5469.599999999999   # M136
